$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill data rows 2-53 with team record values
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 72   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 90   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}

$wb.Save()
